$d = $word.ActiveDocument

# 1. Delete the paragraphs that are fully removed in the target version.
#    These are paragraphs 6-9 (1-indexed) in the original document:
#    "גרוקינג זו תופעה..." through "מאמר חמוד אבל ציפיתי ממנו קצת יותר.."
#    Deleting the range also removes their paragraph marks, merging what
#    remains into a document with the link as paragraph 6.
$p6 = $d.Paragraphs.Item(6)
$p9 = $d.Paragraphs.Item(9)
$rng = $d.Range($p6.Range.Start, $p9.Range.End)
$rng.Delete()

# 2. Update the date in the title line
$d.Content.Find.Execute("01.07.24", $true, $false, $false, $false, $false, $true, 1, $false, "29.06.24", 2)

# 3. Replace the paper title
$d.Content.Find.Execute("Grokfast: Accelerated Grokking by Amplifying Slow Gradients", $true, $false, $false, $false, $false, $true, 1, $false, "What Are the Odds? Language Models Are Capable of Probabilistic Reasoning", 2)

# 4. Replace first body paragraph
$d.Content.Find.Execute("המאמר הזה משך את עיניי משתי סיבות. הסיבה הראשונה היא הופעת מילי Grokking בכותרת. מה זה בעצם Grokking בהקשר של אימון רשתות. אתם בטח יודעים אם אנו מאמנים את הרשת שלנו ליותר מדי זמן (כלומר אפוקים) אז באיזושהי נקודה היא מגיעה למצב של אוורפיט. כלומר הלוס על טריין סט ממשיך לרדת בזמן שהלוס על סט ולידציה מתחיל לעלות כלומר יכולת הכללה של המודל נפגעת. ", $true, $false, $false, $false, $false, $true, 1, $false, "הסקירה הזו הולכת להיות ממש קצרה. לפני ימיים (27.06) סקרתי מאמר שבדק האם מודלי שפה ענקיים מסוגלים לבצע רגרסיה לוגיסטית והגיע למסקנה שבלי עזרה ורמזים מאוד משמעתיים הם לא מצליחים לפתור אותה. ", 2)

# 5. Replace second body paragraph
$d.Content.Find.Execute("אבל אם אנו נמשיך לאמן את הרשת שלנו עוד עוד באיזשהו שלב הלוס על סט ולידציה מתחיל לרדת לאט לאט כלומר יכולת הכללה של המודל משתפרת. כלומר אנו יוצאים מ״משטר האוורפיט״ אחרי שלב מסוים של אימון וזה נקרא grokking. התופעה הזו נחקרת רבות על ידי המדענים בתחום למידה עמוקה אבל אין הבנה מלאה למה זה קורה. השורשים של grokking הזו נמצאים כנראה בתופעה שנקראת double descent.", $true, $false, $false, $false, $false, $true, 1, $false, "הפעם המחברים בדקו האם מודלי שפה מסוגלים ״לנתח התפלגויות הסתברותיות״. למשל אומרים למודל שפה שאיזשהו ערך מפולג גאוסית עם תוחלת 3 ושונות 2 ושואלים אותו מה האחוזון ה-80 של ההתפלגות או לאיזה אחוזון שייכת דגימה בעלת ערך 4. באופן די מפתיע המודל מצליח לא רע בשאלות האלו למרות שקיבל הוראה לא להריץ קוד (זה יכול לעזור כמו שאתם מבינים).", 2)

# 6. Replace third body paragraph
$d.Content.Find.Execute("הסיבה השנייה שבחרתי לסקור את המאמר כי נוכחתה של התמרת פורייה שם אלא אחרי התעמקות קלה התברר שניתן היה להסתדר גם בלעדיו ולהסביר את המאמר בצורה פשוטה יותר בהרבה (מה שאני עושה בסקירה הזו).", $true, $false, $false, $false, $false, $true, 1, $false, "אז מה לדעתכם קורה כאן? איך המודל מצליח לפתור את השאלות האלו?", 2)

# 7. Replace the link text
$d.Content.Find.Execute("https://arxiv.org/abs/2405.20233", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2406.12830", 2)
